$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of row index (1-based, including header row) -> new Estimate value
$updates = @{
    2  = "4"   # U1: 3 -> 4
    3  = "1"   # U2: 3 -> 1
    4  = "4"   # U3: 5 -> 4
    5  = "2"   # U4: 4 -> 2
    7  = "2"   # U6: 1 -> 2
    8  = "2"   # S1: 1 -> 2
    9  = "2"   # S2: 3 -> 2
    10 = "2"   # S3: 4 -> 2
    12 = "4"   # O2: 10 -> 4
    13 = "2"   # O3: 1 -> 2
    14 = "2"   # O4: 5 -> 2
    15 = "2"   # O5: 5 -> 2
    16 = "32"  # Total: 48 -> 32
}

foreach ($rowIndex in $updates.Keys) {
    $cell = $t.Cell($rowIndex, 3)
    $cell.Range.Text = $updates[$rowIndex]
}
